$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '26.088.54'
Set-TextValue 'E2' '  -1.01%  '
Set-TextValue 'D3' '1.668.80'
Set-TextValue 'E3' '  -1.09%  '
Set-TextValue 'E4' '  -0.52%  '
Set-TextValue 'D5' '209.52'
Set-TextValue 'E5' '  -4.05%  '
Set-TextValue 'D6' '0.5249'
Set-TextValue 'E6' '  -3.88%  '
Set-TextValue 'D7' '1.006'
Set-TextValue 'E7' '  -0.47%  '
Set-TextValue 'D8' '0.2656'
Set-TextValue 'E8' '  -2.68%  '
Set-TextValue 'D9' '0.06237'
Set-TextValue 'E9' '  -3.40%  '
Set-TextValue 'D10' '21.09'
Set-TextValue 'E10' '  -4.25%  '
Set-TextValue 'E11' '  -2.08%  '
Set-TextValue 'D12' '1.670.84'
Set-TextValue 'E12' '  -1.21%  '
Set-TextValue 'D13' '4.442'
Set-TextValue 'E13' '  -2.08%  '
Set-TextValue 'D14' '0.5611'
Set-TextValue 'E14' '  -3.46%  '
Set-TextValue 'D15' '0.000008018'
Set-TextValue 'E15' '  -3.79%  '
Set-TextValue 'D16' '65.86'
Set-TextValue 'E16' '  +1.16%  '
Set-TextValue 'D17' '26.151.72'
Set-TextValue 'E17' '  -1.01%  '
Set-TextValue 'D18' '1.006'
Set-TextValue 'E18' '  -0.57%  '
Set-TextValue 'D19' '4.815'
Set-TextValue 'E19' '  -2.56%  '
Set-TextValue 'D20' '10.45'
Set-TextValue 'E20' '  -4.69%  '
Set-TextValue 'D21' '187.08'
Set-TextValue 'E21' '  -2.01%  '
Set-TextValue 'D22' '6.168'
Set-TextValue 'E22' '  -0.99%  '
Set-TextValue 'D23' '1.007'
Set-TextValue 'D24' '146.52'
Set-TextValue 'E24' '  -2.00%  '
Set-TextValue 'D25' '0.1252'
Set-TextValue 'E25' '  -4.23%  '
Set-TextValue 'D26' '7.592'
Set-TextValue 'E26' '  -3.85%  '
Set-TextValue 'D27' '15.79'
Set-TextValue 'E27' '  +0.32%  '
Set-TextValue 'D28' '0.06381'
Set-TextValue 'E28' '  +0.51%  '
Set-TextValue 'D29' '1.341'
Set-TextValue 'E29' '  -5.09%  '
Set-TextValue 'D30' '1.278'
Set-TextValue 'E30' '  -3.79%  '
Set-TextValue 'D31' '3.512'
Set-TextValue 'E31' '  -1.83%  '
Set-TextValue 'D32' '3.457'
Set-TextValue 'E32' '  -3.28%  '
Set-TextValue 'D33' '1.643'
Set-TextValue 'E33' '  -1.84%  '
Set-TextValue 'D34' '1.003'
Set-TextValue 'E34' '  -3.73%  '
Set-TextValue 'B35' 'HuobiToken'
Set-TextValue 'C35' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D35' '2.419'
Set-TextValue 'E35' '  +0.33%  '
Set-TextValue 'B36' 'ImmutableX'
Set-TextValue 'C36' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D36' '0.6037'
Set-TextValue 'E36' '  -2.72%  '
Set-TextValue 'D37' '2.702'
Set-TextValue 'E37' '  -0.64%  '
Set-TextValue 'D38' '6.125'
Set-TextValue 'E38' '  -1.83%  '
Set-TextValue 'D39' '1.098.42'
Set-TextValue 'E39' '  -1.58%  '
Set-TextValue 'D40' '0.01607'
Set-TextValue 'E40' '  -1.49%  '
Set-TextValue 'D41' '0.8631'
Set-TextValue 'E41' '  -1.84%  '
Set-TextValue 'D42' '1.005'
Set-TextValue 'E42' '  -0.99%  '
Set-TextValue 'D43' '99.76'
Set-TextValue 'E43' '  -1.17%  '
Set-TextValue 'D44' '1.827.06'
Set-TextValue 'E44' '  -0.69%  '
Set-TextValue 'E45' '  -1.07%  '
Set-TextValue 'D46' '56.59'
Set-TextValue 'E46' '  -1.26%  '
Set-TextValue 'D47' '1.003'
Set-TextValue 'E47' '  -0.98%  '
Set-TextValue 'D48' '0.05257'
Set-TextValue 'E48' '  -0.22%  '
Set-TextValue 'D49' '7.946'
Set-TextValue 'E49' '  -3.27%  '
Set-TextValue 'D50' '0.4267'
Set-TextValue 'E50' '  -0.91%  '
Set-TextValue 'D51' '5.915'
Set-TextValue 'E51' '  -2.02%  '
